# "check same action with 1 row" - replicate row 7 ("**" action) into a block
# of rows (7-18) each carrying the 2017/12/25 date + duplicated 89/89 pair and
# the D/E amount pairs, then push the trailing "**" marker row down to row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

$rows = @(
    @(89, "2017/12/25", 89, 25000,  25029),
    @(89, "2017/12/25", 89, 35000,  35029),
    @(89, "2017/12/25", 89, 350000, 350029),
    @(89, "2017/12/25", 89, 150000, 150029),
    @(89, "2017/12/25", 89, 34000,  34029),
    @(89, "2017/12/25", 89, 29000,  29029),
    @(89, "2017/12/25", 89, 29000,  29002),
    @(89, "2017/12/25", 89, 29000,  29002),
    @(89, "2017/12/25", 89, 29000,  29002),
    @(89, "2017/12/25", 89, 25000,  25002),
    @(89, "2017/12/25", 89, 35000,  35002),
    @(89, "2017/12/25", 89, 34000,  34002)
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# trailing marker row, moved from row 7 down to row 19
$ws.Cells.Item(19, 1).Value = "**"
